$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A67").Value = "22-11-2025"
$ws.Range("B67").Value = "The price of gold in India today is ₹12,584 per gram for 24 karat gold, ₹11,535 per gram for 22 karat gold and ₹9,438 per gram for 18 karat gold (also called 999 gold)."
